$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new data rows
$ws.Range("A31").Value = 901
$ws.Range("A32").Value = 1109

# Update the view: scroll so row 7 is the top row, and select column B (whole column)
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("B1:B1048576").Select()
